# Fruta / hortaliza, semanal
# Insert two new weekly records (rows 54 and 55) for
# "Femacal de La Calera - Papaya", pushing the existing historical
# rows down by two positions (old row 54 -> new row 56, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 54; this shifts every
# existing row at/after 54 down by two rows (54->56, 55->57, ...).
$ws.Range("A54:A55").EntireRow.Insert()

# --- Row 54 : "Primera" quality -----------------------------------
$ws.Cells.Item(54, 1).Value2 = 3
$ws.Cells.Item(54, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(54, 3).Value = "Coquimbo"
$ws.Cells.Item(54, 4).Value2 = 45219
$ws.Cells.Item(54, 5).Value2 = 5
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value2 = 100108
$ws.Cells.Item(54, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(54, 9).Value2 = 100108004
$ws.Cells.Item(54, 10).Value = "Papaya"
$ws.Cells.Item(54, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(54, 12).Value = "Primera"
$ws.Cells.Item(54, 13).Value2 = 40
$ws.Cells.Item(54, 14).Value2 = 16000
$ws.Cells.Item(54, 15).Value2 = 16000
$ws.Cells.Item(54, 16).Value2 = 16000
$ws.Cells.Item(54, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(54, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(54, 19).Value2 = 1600
$ws.Cells.Item(54, 20).Value2 = 10

# --- Row 55 : "Segunda" quality ------------------------------------
$ws.Cells.Item(55, 1).Value2 = 3
$ws.Cells.Item(55, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(55, 3).Value = "Coquimbo"
$ws.Cells.Item(55, 4).Value2 = 45219
$ws.Cells.Item(55, 5).Value2 = 5
$ws.Cells.Item(55, 6).Value = "Fruta"
$ws.Cells.Item(55, 7).Value2 = 100108
$ws.Cells.Item(55, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(55, 9).Value2 = 100108004
$ws.Cells.Item(55, 10).Value = "Papaya"
$ws.Cells.Item(55, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(55, 12).Value = "Segunda"
$ws.Cells.Item(55, 13).Value2 = 30
$ws.Cells.Item(55, 14).Value2 = 13000
$ws.Cells.Item(55, 15).Value2 = 13000
$ws.Cells.Item(55, 16).Value2 = 13000
$ws.Cells.Item(55, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(55, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(55, 19).Value2 = 1300
$ws.Cells.Item(55, 20).Value2 = 10
